# Apply the edit described in the diff to the "incidence2018_plus" sheet:
#  - Column C formulas for rows 27..122 change divisor from /2 to /8
#  - Selection changes from F16 to C27:C122 (active cell C27)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("incidence2018_plus")
$ws.Activate()

for ($row = 27; $row -le 122; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Formula = "=prevalence2018!C$row/8"
}

$rangeToSelect = $ws.Range("C27:C122")
$rangeToSelect.Select()
